$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 12.93024956307966
$ws.Range("D2").Value2 = 4.664589018196449
$ws.Range("E2").Value2 = 18.89347335491468
$ws.Range("F2").Value2 = 21.92062565984711
$ws.Range("G2").Value2 = 24.64875663200993
$ws.Range("H2").Value2 = 13.32483823795075
$ws.Range("I2").Value2 = 26.27469423851406
$ws.Range("K2").Value2 = 11.59156521949556
$ws.Range("L2").Value2 = 8.687228697921983
$ws.Range("M2").Value2 = 13.7878532964395
$ws.Range("O2").Value2 = 19.72122102775513
$ws.Range("B3").Value2 = 12.79774399502409
$ws.Range("D3").Value2 = 4.579613810151755
$ws.Range("E3").Value2 = 18.96171057244061
$ws.Range("F3").Value2 = 21.95874607700489
$ws.Range("G3").Value2 = 24.71339003700496
$ws.Range("H3").Value2 = 13.36960349124182
$ws.Range("I3").Value2 = 26.40554002072659
$ws.Range("K3").Value2 = 11.32941378140598
$ws.Range("L3").Value2 = 8.667293570281389
$ws.Range("M3").Value2 = 13.75292521618595
$ws.Range("O3").Value2 = 19.79337512190975
$ws.Range("B4").Value2 = 12.71750729189402
$ws.Range("D4").Value2 = 4.52594774618751
$ws.Range("E4").Value2 = 19.00592029064736
$ws.Range("F4").Value2 = 21.98814564910303
$ws.Range("G4").Value2 = 24.76159673066486
$ws.Range("H4").Value2 = 13.39915720482396
$ws.Range("I4").Value2 = 26.49022742279231
$ws.Range("K4").Value2 = 11.16404145506428
$ws.Range("L4").Value2 = 8.656186417279226
$ws.Range("M4").Value2 = 13.73306753378433
$ws.Range("O4").Value2 = 19.84193455135058
$ws.Range("B5").Value2 = 12.68512420312044
$ws.Range("D5").Value2 = 4.503718719450474
$ws.Range("E5").Value2 = 19.02451905608566
$ws.Range("F5").Value2 = 22.00163050839801
$ws.Range("G5").Value2 = 24.78337536232264
$ws.Range("H5").Value2 = 13.41172068928244
$ws.Range("I5").Value2 = 26.52583382373307
$ws.Range("K5").Value2 = 11.09560384751504
$ws.Range("L5").Value2 = 8.651948214560967
$ws.Range("M5").Value2 = 13.725380101747
$ws.Range("O5").Value2 = 19.8627915475543
$ws.Range("B6").Value2 = 12.67976690126858
$ws.Range("D6").Value2 = 4.500006342303956
$ws.Range("E6").Value2 = 19.02764262337566
$ws.Range("F6").Value2 = 22.00396042062719
$ws.Range("G6").Value2 = 24.78712030132407
$ws.Range("H6").Value2 = 13.41383826466537
$ws.Range("I6").Value2 = 26.53181248079676
$ws.Range("K6").Value2 = 11.08417847748703
$ws.Range("L6").Value2 = 8.651261949212454
$ws.Range("M6").Value2 = 13.72412820497019
$ws.Range("O6").Value2 = 19.86631932800858
$ws.Range("B7").Value2 = 12.71706924913647
$ws.Range("D7").Value2 = 4.525649393412679
$ws.Range("E7").Value2 = 19.00616875759887
$ws.Range("F7").Value2 = 21.98832142396421
$ws.Range("G7").Value2 = 24.76188181668681
$ws.Range("H7").Value2 = 13.39932453424428
$ws.Range("I7").Value2 = 26.49070318397438
$ws.Range("K7").Value2 = 11.1631226364834
$ws.Range("L7").Value2 = 8.656128089103115
$ws.Range("M7").Value2 = 13.73296221278554
$ws.Range("O7").Value2 = 19.84221151152418
$ws.Range("B8").Value2 = 12.88434795779461
$ws.Range("D8").Value2 = 4.635607932415271
$ws.Range("E8").Value2 = 18.91652271656405
$ws.Range("F8").Value2 = 21.93252411848104
$ws.Range("G8").Value2 = 24.66926858397772
$ws.Range("H8").Value2 = 13.33984425880334
$ws.Range("I8").Value2 = 26.31890895438804
$ws.Range("K8").Value2 = 11.50212679781438
$ws.Range("L8").Value2 = 8.68012182570348
$ws.Range("M8").Value2 = 13.77548373046485
$ws.Range("O8").Value2 = 19.74521531082531
$ws.Range("B9").Value2 = 13.21986808277146
$ws.Range("D9").Value2 = 4.838774731061203
$ws.Range("E9").Value2 = 18.75899587883315
$ws.Range("F9").Value2 = 21.87076171747729
$ws.Range("G9").Value2 = 24.55561833817346
$ws.Range("H9").Value2 = 13.23960168853439
$ws.Range("I9").Value2 = 26.01640825931438
$ws.Range("K9").Value2 = 12.12938668181021
$ws.Range("L9").Value2 = 8.736018935819398
$ws.Range("M9").Value2 = 13.87122527646213
$ws.Range("O9").Value2 = 19.58885385445796
$ws.Range("B10").Value2 = 13.469025554495
$ws.Range("D10").Value2 = 4.979653613540584
$ws.Range("E10").Value2 = 18.65429416797067
$ws.Range("F10").Value2 = 21.85453622047907
$ws.Range("G10").Value2 = 24.51396373990518
$ws.Range("H10").Value2 = 13.175939077811
$ws.Range("I10").Value2 = 25.81497135463927
$ws.Range("K10").Value2 = 12.56423624698955
$ws.Range("L10").Value2 = 8.782273343948292
$ws.Range("M10").Value2 = 13.94875560977929
$ws.Range("O10").Value2 = 19.49470807059359
$ws.Range("B11").Value2 = 13.58254293148735
$ws.Range("D11").Value2 = 5.041769839630732
$ws.Range("E11").Value2 = 18.60903636715469
$ws.Range("F11").Value2 = 21.8534918700582
$ws.Range("G11").Value2 = 24.50416501299729
$ws.Range("H11").Value2 = 13.14914334838892
$ws.Range("I11").Value2 = 25.72781904011353
$ws.Range("K11").Value2 = 12.75583128938513
$ws.Range("L11").Value2 = 8.804394206953599
$ws.Range("M11").Value2 = 13.98550846865672
$ws.Range("O11").Value2 = 19.45640058874822
$ws.Range("B12").Value2 = 13.625520368463
$ws.Range("D12").Value2 = 5.064996611390225
$ws.Range("E12").Value2 = 18.59223777807408
$ws.Range("F12").Value2 = 21.85400694914286
$ws.Range("G12").Value2 = 24.50177366342345
$ws.Range("H12").Value2 = 13.13930765960868
$ws.Range("I12").Value2 = 25.69545894756646
$ws.Range("K12").Value2 = 12.82744507837354
$ws.Range("L12").Value2 = 8.812921681026745
$ws.Range("M12").Value2 = 13.99963210901612
$ws.Range("O12").Value2 = 19.44254606107686
$ws.Range("B13").Value2 = 13.61626535362571
$ws.Range("D13").Value2 = 5.060007626006639
$ws.Range("E13").Value2 = 18.59584057606478
$ws.Range("F13").Value2 = 21.85385553908128
$ws.Range("G13").Value2 = 24.50222996958009
$ws.Range("H13").Value2 = 13.14141210738149
$ws.Range("I13").Value2 = 25.70239972356326
$ws.Range("K13").Value2 = 12.81206420573474
$ws.Range("L13").Value2 = 8.811078503415626
$ws.Range("M13").Value2 = 13.99658128064699
$ws.Range("O13").Value2 = 19.44550087103587
$ws.Range("B14").Value2 = 13.5860790881716
$ws.Range("D14").Value2 = 5.043686704022219
$ws.Range("E14").Value2 = 18.60764754022036
$ws.Range("F14").Value2 = 21.85351600354716
$ws.Range("G14").Value2 = 24.5039418209805
$ws.Range("H14").Value2 = 13.14832792203997
$ws.Range("I14").Value2 = 25.72514388828949
$ws.Range("K14").Value2 = 12.76174206236915
$ws.Range("L14").Value2 = 8.805092771538796
$ws.Range("M14").Value2 = 13.98666634027291
$ws.Range("O14").Value2 = 19.45524769809154
$ws.Range("B15").Value2 = 13.56758693344931
$ws.Range("D15").Value2 = 5.033650870193103
$ws.Range("E15").Value2 = 18.61492382388394
$ws.Range("F15").Value2 = 21.85342657662041
$ws.Range("G15").Value2 = 24.50516226069878
$ws.Range("H15").Value2 = 13.15260459480618
$ws.Range("I15").Value2 = 25.73915897362757
$ws.Range("K15").Value2 = 12.73079473672753
$ws.Range("L15").Value2 = 8.801445839442339
$ws.Range("M15").Value2 = 13.98061977795631
$ws.Range("O15").Value2 = 19.46130282950209
$ws.Range("B16").Value2 = 13.46160750085929
$ws.Range("D16").Value2 = 4.975553483535179
$ws.Range("E16").Value2 = 18.65729946738486
$ws.Range("F16").Value2 = 21.85473194275591
$ws.Range("G16").Value2 = 24.51478855645035
$ws.Range("H16").Value2 = 13.17773380030892
$ws.Range("I16").Value2 = 25.82075699130233
$ws.Range("K16").Value2 = 12.55158599374293
$ws.Range("L16").Value2 = 8.780849039614994
$ws.Range("M16").Value2 = 13.94638296943444
$ws.Range("O16").Value2 = 19.49730265619909
$ws.Range("B17").Value2 = 13.39661272979247
$ws.Range("D17").Value2 = 4.939399212426568
$ws.Range("E17").Value2 = 18.68390191124465
$ws.Range("F17").Value2 = 21.85715549389888
$ws.Range("G17").Value2 = 24.523040391449
$ws.Range("H17").Value2 = 13.19370418002674
$ws.Range("I17").Value2 = 25.87196136240808
$ws.Range("K17").Value2 = 12.44002115306293
$ws.Range("L17").Value2 = 8.768486947147176
$ws.Range("M17").Value2 = 13.92575477700349
$ws.Range("O17").Value2 = 19.52054629486856
$ws.Range("B18").Value2 = 13.35924789207137
$ws.Range("D18").Value2 = 4.918419444572051
$ws.Range("E18").Value2 = 18.69942623998398
$ws.Range("F18").Value2 = 21.85914598407342
$ws.Range("G18").Value2 = 24.52864775781297
$ws.Range("H18").Value2 = 13.20309367333375
$ws.Range("I18").Value2 = 25.9018347763935
$ws.Range("K18").Value2 = 12.3752689622582
$ws.Range("L18").Value2 = 8.761478483705369
$ws.Range("M18").Value2 = 13.91403010958697
$ws.Range("O18").Value2 = 19.53434075779291
$ws.Range("B19").Value2 = 13.34660101650947
$ws.Range("D19").Value2 = 4.911284696004016
$ws.Range("E19").Value2 = 18.70472090896199
$ws.Range("F19").Value2 = 21.85992239190136
$ws.Range("G19").Value2 = 24.53069409039991
$ws.Range("H19").Value2 = 13.20630778963508
$ws.Range("I19").Value2 = 25.9120219426949
$ws.Range("K19").Value2 = 12.35324626527414
$ws.Range("L19").Value2 = 8.759123169944843
$ws.Range("M19").Value2 = 13.91008462217561
$ws.Range("O19").Value2 = 19.53908432397959
$ws.Range("B20").Value2 = 13.40352987547948
$ws.Range("D20").Value2 = 4.943267110323806
$ws.Range("E20").Value2 = 18.68104693495101
$ws.Range("F20").Value2 = 21.85683576676891
$ws.Range("G20").Value2 = 24.52207281704104
$ws.Range("H20").Value2 = 13.19198301822398
$ws.Range("I20").Value2 = 25.86646690712913
$ws.Range("K20").Value2 = 12.45195809885128
$ws.Range("L20").Value2 = 8.769792399929175
$ws.Range("M20").Value2 = 13.92793623288635
$ws.Range("O20").Value2 = 19.51802793632565
$ws.Range("B21").Value2 = 13.59494604390113
$ws.Range("D21").Value2 = 5.048488660122858
$ws.Range("E21").Value2 = 18.60417034475197
$ws.Range("F21").Value2 = 21.85359103021531
$ws.Range("G21").Value2 = 24.5034031844587
$ws.Range("H21").Value2 = 13.1462881318573
$ws.Range("I21").Value2 = 25.71844595128253
$ws.Range("K21").Value2 = 12.77654871996638
$ws.Range("L21").Value2 = 8.806846870164819
$ws.Range("M21").Value2 = 13.98957306325654
$ws.Range("O21").Value2 = 19.45236711883384
$ws.Range("B22").Value2 = 13.71997888641063
$ws.Range("D22").Value2 = 5.115530824802033
$ws.Range("E22").Value2 = 18.55590558184269
$ws.Range("F22").Value2 = 21.85677696393389
$ws.Range("G22").Value2 = 24.49889159782144
$ws.Range("H22").Value2 = 13.11823821313239
$ws.Range("I22").Value2 = 25.62545036909769
$ws.Range("K22").Value2 = 12.98319555388933
$ws.Range("L22").Value2 = 8.831940807512973
$ws.Range("M22").Value2 = 14.03105425487012
$ws.Range("O22").Value2 = 19.41325323645947
$ws.Range("B23").Value2 = 13.65326395092094
$ws.Range("D23").Value2 = 5.079910857641599
$ws.Range("E23").Value2 = 18.58148483887982
$ws.Range("F23").Value2 = 21.85459143271386
$ws.Range("G23").Value2 = 24.50059506381364
$ws.Range("H23").Value2 = 13.13304297392479
$ws.Range("I23").Value2 = 25.67474186079151
$ws.Range("K23").Value2 = 12.87342056011261
$ws.Range("L23").Value2 = 8.818468995917367
$ws.Range("M23").Value2 = 14.00880780487093
$ws.Range("O23").Value2 = 19.43378086999845
$ws.Range("B24").Value2 = 13.40040262656361
$ws.Range("D24").Value2 = 4.941519037647947
$ws.Range("E24").Value2 = 18.68233695249139
$ws.Range("F24").Value2 = 21.85697845527356
$ws.Range("G24").Value2 = 24.52250756868029
$ws.Range("H24").Value2 = 13.19276050785005
$ws.Range("I24").Value2 = 25.86894959416681
$ws.Range("K24").Value2 = 12.4465633087434
$ws.Range("L24").Value2 = 8.769201897052614
$ws.Range("M24").Value2 = 13.92694957638879
$ws.Range("O24").Value2 = 19.5191651426213
$ws.Range("B25").Value2 = 13.12849675972332
$ws.Range("D25").Value2 = 4.785231474155397
$ws.Range("E25").Value2 = 18.79966605329109
$ws.Range("F25").Value2 = 21.88235526883103
$ws.Range("G25").Value2 = 24.5790409122325
$ws.Range("H25").Value2 = 13.26496559905321
$ws.Range("I25").Value2 = 26.09457727709473
$ws.Range("K25").Value2 = 11.96405083400017
$ws.Range("L25").Value2 = 8.71997075927777
$ws.Range("M25").Value2 = 13.84403617928549
$ws.Range("O25").Value2 = 19.6275197161305
